$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.973.40'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '2.233.87'
$ws.Range('E3').Value = '  -2.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.11'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.56%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '267.73'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.44%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.96%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.602'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '46.12'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0927'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.65%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.21'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('E13').Value = '  -2.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.29'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.876'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.88%  '
$ws.Range('D16').Value = '2.571.46'
$ws.Range('E16').Value = '  -1.67%  '
$ws.Range('D17').Value = '2.250.50'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D18').Value = '43.003.32'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('E19').Value = '  -1.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.70'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.62'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.35'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -6.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.81'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.92'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.25'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.21'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.99%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.01%  '
$ws.Range('E29').Value = '  -1.18%  '
$ws.Range('E30').Value = '  -1.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.23'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.09'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0897'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.57'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.92%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.21'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +7.19%  '
$ws.Range('E37').Value = '  +0.81%  '
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.106'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.51'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.28'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -6.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '70.91'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -8.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.231'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.32%  '
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.67'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -9.73%  '
$ws.Range('E46').Value = '  -4.55%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.25'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.648'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.87%  '
$ws.Range('B49').Value = 'FraxShare'
$ws.Range('C49').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.40'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0988'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '99.86'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -4.64%  '
